# Cham cong function - update last_edited_time and derived metrics
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Thang 7)
$ws.Range("D7").Value = "2024-07-06T13:10:00.000Z"
$ws.Range("W7").Value = 31575000
$ws.Range("AA7").Value = 46455000

# Row 8 (Thang 6)
$ws.Range("D8").Value = "2024-07-06T13:10:00.000Z"
$ws.Range("AK8").Value = 103

# Row 9 (Thang 5)
$ws.Range("D9").Value = "2024-07-06T13:10:00.000Z"

# Row 10 (Thang 4)
$ws.Range("D10").Value = "2024-07-06T13:10:00.000Z"

# Row 11 (Thang 3)
$ws.Range("D11").Value = "2024-07-06T13:10:00.000Z"

# Row 12 (Thang 2)
$ws.Range("D12").Value = "2024-07-06T13:10:00.000Z"
